$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.969.80"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
$ws.Range("D3").Value = "3.216.75"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'578.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "

# Row 6
$ws.Range("D6").Value = "'142.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.38%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "3.205.14"
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").Value = "'0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.43%  "

# Row 11
$ws.Range("D11").Value = "'6.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.90%  "

# Row 12
$ws.Range("D12").Value = "'0.479"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "

# Row 13
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.72%  "

# Row 14
$ws.Range("D14").Value = "'35.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.78%  "

# Row 15
$ws.Range("D15").Value = "3.730.93"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").Value = "66.993.52"
$ws.Range("E16").Value = "  +0.68%  "

# Row 17
$ws.Range("D17").Value = "3.225.96"
$ws.Range("E17").Value = "  +0.25%  "

# Row 18
$ws.Range("E18").Value = "  -1.87%  "

# Row 19
$ws.Range("D19").Value = "'6.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.19%  "

# Row 20
$ws.Range("D20").Value = "'499.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "

# Row 21
$ws.Range("D21").Value = "'14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("D22").Value = "'0.715"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.67%  "

# Row 23
$ws.Range("D23").Value = "'7.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.81%  "

# Row 24
$ws.Range("D24").Value = "'81.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "

# Row 25
$ws.Range("D25").Value = "'12.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "

# Row 26
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("D27").Value = "'3.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.34%  "

# Row 28
$ws.Range("D28").Value = "'7.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.79%  "

# Row 29
$ws.Range("D29").Value = "'2.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.58%  "

# Row 30
$ws.Range("D30").Value = "'27.80"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.03%  "

# Row 32
$ws.Range("E32").Value = "  -1.98%  "

# Row 33
$ws.Range("E33").Value = "  -0.07%  "

# Row 34
$ws.Range("D34").Value = "'509.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.05%  "

# Row 35
$ws.Range("D35").Value = "'6.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.43%  "

# Row 36
$ws.Range("D36").Value = "'54.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.52%  "

# Row 37
$ws.Range("D37").Value = "'5.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.90%  "

# Row 38
$ws.Range("D38").Value = "'0.0415"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.58%  "

# Row 39
$ws.Range("D39").Value = "'0.0814"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.56%  "

# Row 40
$ws.Range("D40").Value = "'8.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.84%  "

# Row 41
$ws.Range("D41").Value = "'0.118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.82%  "

# Row 42
$ws.Range("D42").Value = "2.867.58"
$ws.Range("E42").Value = "  -0.40%  "

# Row 43
$ws.Range("D43").Value = "'2.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.24%  "

# Row 44
$ws.Range("D44").Value = "'0.251"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "

# Row 46
$ws.Range("D46").Value = "'123.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.07%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'24.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.12%  "

# Row 49
$ws.Range("D49").Value = "0.0₃0525"
$ws.Range("E49").Value = "  -8.23%  "

# Row 50
$ws.Range("E50").Value = "  -2.12%  "

# Row 51
$ws.Range("E51").Value = "  -12.06%  "
